$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-10 from 45233 to 45243
$ws.Range("C2:C10").Value = 45243
